$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force plain-numeric-looking strings to be stored as text,
# matching the original workbook where these price cells are inline strings, not numbers.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$ws.Range('D2').Value = '42.300.46'
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').Value = '2.243.57'
$ws.Range('E3').Value = '  -4.50%  '
$ws.Range('E4').Value = '  +0.16%  '
$scratch.Value = '232.86'
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -3.04%  '
$scratch.Value = '0.633'
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -5.53%  '
$scratch.Value = '69.52'
$scratch.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -3.80%  '
$ws.Range('E8').Value = '  +0.12%  '
$scratch.Value = '0.559'
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -5.87%  '
$scratch.Value = '0.0993'
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.95%  '
$scratch.Value = '58.33'
$scratch.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +0.26%  '
$scratch.Value = '35.85'
$scratch.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +9.12%  '
$scratch.Value = '0.105'
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  -3.15%  '
$scratch.Value = '6.77'
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  -6.55%  '
$ws.Range('D15').Value = '2.577.75'
$ws.Range('E15').Value = '  -4.46%  '
$scratch.Value = '15.03'
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  -7.86%  '
$scratch.Value = '0.860'
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -4.75%  '
$ws.Range('D18').Value = '2.240.02'
$ws.Range('E18').Value = '  -4.87%  '
$ws.Range('D19').Value = '42.108.11'
$ws.Range('E19').Value = '  -3.84%  '
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('E20').Value = '  -5.03%  '
$scratch.Value = '6.24'
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -7.31%  '
$scratch.Value = '73.38'
$scratch.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -6.17%  '
$scratch.Value = '236.41'
$scratch.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -7.13%  '
$scratch.Value = '2.01'
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +4.33%  '
$scratch.Value = '0.999'
$scratch.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -2.14%  '
$scratch.Value = '2.35'
$scratch.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -5.48%  '
$scratch.Value = '10.03'
$scratch.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -4.12%  '
$scratch.Value = '2.19'
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -3.96%  '
$scratch.Value = '169.04'
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -4.41%  '
$scratch.Value = '20.62'
$scratch.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -7.89%  '
$ws.Range('E32').Value = '  -6.29%  '
$scratch.Value = '0.127'
$scratch.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -6.85%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$scratch.Value = '0.0715'
$scratch.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$scratch.Value = '5.32'
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -1.24%  '
$scratch.Value = '4.73'
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -7.78%  '
$scratch.Value = '3.62'
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -3.54%  '
$scratch.Value = '21.92'
$scratch.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +16.54%  '
$scratch.Value = '2.25'
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -5.28%  '
$scratch.Value = '6.04'
$scratch.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -5.97%  '
$ws.Range('E41').Value = '  -3.47%  '
$scratch.Value = '66.95'
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +0.19%  '
$scratch.Value = '4.97'
$scratch.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -4.06%  '
$scratch.Value = '8.98'
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -2.20%  '
$scratch.Value = '0.192'
$scratch.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -3.47%  '
$ws.Range('E46').Value = '  -7.47%  '
$ws.Range('E47').Value = '  +0.01%  '
$scratch.Value = '4.37'
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +6.84%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$scratch.Value = '1.18'
$scratch.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -4.73%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$scratch.Value = '2.36'
$scratch.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -5.00%  '
$scratch.Value = '10.07'
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +7.29%  '

$scratch.Clear()
$excel.CutCopyMode = 0

